# Add market value calculation (holding_value column) to the ETF top
# holdings sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell E1, matching formatting of the existing header row (D1).
$ws.Range("E1").Value = "holding_value"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats

# Market value for each of the 70 holdings (rows 2..71).
$holdingValues = @(
    365.02, 301.33, 251.34, 201.05, 165.38, 154.55, 149.64, 143.86, 126.33, 122.96,
    444.17, 395.26, 390.75, 380.14, 366.54, 361.72, 361.04, 355.37, 353.76, 350.45,
    410.14, 338.55, 282.36, 205.47, 177.88, 173.63, 152.99, 138.97, 138.26, 106.03,
    1353.17, 931.77, 572.95, 558.55, 458.65, 455.83, 454.65, 359.48, 256.74, 239.97,
    583.08, 462.98, 353.28, 139.79, 75.90000000000001, 61.23, 52.46, 50.27, 49.06, 46.38,
    1863.97, 1538.78, 1283.49, 934.28, 789.23, 631.76, 627.89, 626.26, 485.38, 353.64,
    910.49, 791.22, 659.97, 475.09, 405.7, 322.83, 322, 321.14, 251.06, 182.07
)

for ($i = 0; $i -lt $holdingValues.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $holdingValues[$i]
}
